# Update row 29 (2025Q3) metrics in metricas_recorrencia_trimestral
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 82
$ws.Range("D29").Value = 14
$ws.Range("E29").Value = 68
$ws.Range("F29").Value = 2.409638554216868
